$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.6869109999999999
$ws.Cells.Item(2, 8).Value = 2.060733
$ws.Cells.Item(2, 9).Value = 0.08184737208886859
$ws.Cells.Item(2, 10).Value = 0.08184737208886859
$ws.Cells.Item(2, 13).Value = 10.43365533333333
$ws.Cells.Item(2, 14).Value = 31.300966
$ws.Cells.Item(2, 15).Value = 0.1673546557507891
$ws.Cells.Item(2, 16).Value = 0.1673546557507891
$ws.Cells.Item(2, 17).Value = 7.166992618675334
$ws.Cells.Item(2, 18).Value = 64.50293356807801
$ws.Cells.Item(2, 19).Value = 0.01369753878003935
$ws.Cells.Item(2, 20).Value = 0.01369753878003935

$ws.Cells.Item(3, 7).Value = 0.6869109999999999
$ws.Cells.Item(3, 8).Value = 2.060733
$ws.Cells.Item(3, 9).Value = 0.08184737208886859
$ws.Cells.Item(3, 10).Value = 0.08184737208886859
$ws.Cells.Item(3, 15).Value = 0.5360787938719054
$ws.Cells.Item(3, 16).Value = 0.5360787938719054
$ws.Cells.Item(3, 17).Value = 22.957668799066
$ws.Cells.Item(3, 18).Value = 206.619019191594
$ws.Cells.Item(3, 19).Value = 0.04387664051098573
$ws.Cells.Item(3, 20).Value = 0.04387664051098573

$ws.Cells.Item(4, 7).Value = 0.6869109999999999
$ws.Cells.Item(4, 8).Value = 2.060733
$ws.Cells.Item(4, 9).Value = 0.08184737208886859
$ws.Cells.Item(4, 10).Value = 0.08184737208886859
$ws.Cells.Item(4, 13).Value = 18.162159
$ws.Cells.Item(4, 14).Value = 54.486477
$ws.Cells.Item(4, 15).Value = 0.2913189836188534
$ws.Cells.Item(4, 16).Value = 0.2913189836188534
$ws.Cells.Item(4, 17).Value = 12.475786800849
$ws.Cells.Item(4, 18).Value = 112.282081207641
$ws.Cells.Item(4, 19).Value = 0.02384369324880331
$ws.Cells.Item(4, 20).Value = 0.02384369324880331

$ws.Cells.Item(5, 7).Value = 0.6869109999999999
$ws.Cells.Item(5, 8).Value = 2.060733
$ws.Cells.Item(5, 9).Value = 0.08184737208886859
$ws.Cells.Item(5, 10).Value = 0.08184737208886859
$ws.Cells.Item(5, 13).Value = 0.3271573333333334
$ws.Cells.Item(5, 14).Value = 0.981472
$ws.Cells.Item(5, 15).Value = 0.005247566758452071
$ws.Cells.Item(5, 16).Value = 0.00524756675845207
$ws.Cells.Item(5, 17).Value = 0.2247279709973333
$ws.Cells.Item(5, 18).Value = 2.022551738976
$ws.Cells.Item(5, 19).Value = 0.0004294995490402046
$ws.Cells.Item(5, 20).Value = 0.0004294995490402046

$ws.Cells.Item(6, 9).Value = 0.1692706982521157
$ws.Cells.Item(6, 10).Value = 0.1692706982521157
$ws.Cells.Item(6, 13).Value = 10.43365533333333
$ws.Cells.Item(6, 14).Value = 31.300966
$ws.Cells.Item(6, 15).Value = 0.1673546557507891
$ws.Cells.Item(6, 16).Value = 0.1673546557507891
$ws.Cells.Item(6, 17).Value = 14.82224552809956
$ws.Cells.Item(6, 18).Value = 133.400209752896
$ws.Cells.Item(6, 19).Value = 0.02832823943467852
$ws.Cells.Item(6, 20).Value = 0.02832823943467852

$ws.Cells.Item(7, 9).Value = 0.1692706982521157
$ws.Cells.Item(7, 10).Value = 0.1692706982521157
$ws.Cells.Item(7, 15).Value = 0.5360787938719054
$ws.Cells.Item(7, 16).Value = 0.5360787938719054
$ws.Cells.Item(7, 19).Value = 0.09074243175684943
$ws.Cells.Item(7, 20).Value = 0.09074243175684943

$ws.Cells.Item(8, 9).Value = 0.1692706982521157
$ws.Cells.Item(8, 10).Value = 0.1692706982521157
$ws.Cells.Item(8, 13).Value = 18.162159
$ws.Cells.Item(8, 14).Value = 54.486477
$ws.Cells.Item(8, 15).Value = 0.2913189836188534
$ws.Cells.Item(8, 16).Value = 0.2913189836188534
$ws.Cells.Item(8, 17).Value = 25.801502102368
$ws.Cells.Item(8, 18).Value = 232.213518921312
$ws.Cells.Item(8, 19).Value = 0.04931176777125997
$ws.Cells.Item(8, 20).Value = 0.04931176777125997

$ws.Cells.Item(9, 9).Value = 0.1692706982521157
$ws.Cells.Item(9, 10).Value = 0.1692706982521157
$ws.Cells.Item(9, 13).Value = 0.3271573333333334
$ws.Cells.Item(9, 14).Value = 0.981472
$ws.Cells.Item(9, 15).Value = 0.005247566758452071
$ws.Cells.Item(9, 16).Value = 0.00524756675845207
$ws.Cells.Item(9, 17).Value = 0.4647658146702223
$ws.Cells.Item(9, 18).Value = 4.182892332032
$ws.Cells.Item(9, 19).Value = 0.0008882592893277734
$ws.Cells.Item(9, 20).Value = 0.0008882592893277733

$ws.Cells.Item(10, 7).Value = 2.651481
$ws.Cells.Item(10, 8).Value = 7.954443000000001
$ws.Cells.Item(10, 9).Value = 0.3159313972167653
$ws.Cells.Item(10, 10).Value = 0.3159313972167653
$ws.Cells.Item(10, 13).Value = 10.43365533333333
$ws.Cells.Item(10, 14).Value = 31.300966
$ws.Cells.Item(10, 15).Value = 0.1673546557507891
$ws.Cells.Item(10, 16).Value = 0.1673546557507891
$ws.Cells.Item(10, 17).Value = 27.66463887688201
$ws.Cells.Item(10, 18).Value = 248.9817498919381
$ws.Cells.Item(10, 19).Value = 0.05287259022207756
$ws.Cells.Item(10, 20).Value = 0.05287259022207756

$ws.Cells.Item(11, 7).Value = 2.651481
$ws.Cells.Item(11, 8).Value = 7.954443000000001
$ws.Cells.Item(11, 9).Value = 0.3159313972167653
$ws.Cells.Item(11, 10).Value = 0.3159313972167653
$ws.Cells.Item(11, 15).Value = 0.5360787938719054
$ws.Cells.Item(11, 16).Value = 0.5360787938719054
$ws.Cells.Item(11, 17).Value = 88.616753298486
$ws.Cells.Item(11, 18).Value = 797.5507796863741
$ws.Cells.Item(11, 19).Value = 0.1693641223662294
$ws.Cells.Item(11, 20).Value = 0.1693641223662294

$ws.Cells.Item(12, 7).Value = 2.651481
$ws.Cells.Item(12, 8).Value = 7.954443000000001
$ws.Cells.Item(12, 9).Value = 0.3159313972167653
$ws.Cells.Item(12, 10).Value = 0.3159313972167653
$ws.Cells.Item(12, 13).Value = 18.162159
$ws.Cells.Item(12, 14).Value = 54.486477
$ws.Cells.Item(12, 15).Value = 0.2913189836188534
$ws.Cells.Item(12, 16).Value = 0.2913189836188534
$ws.Cells.Item(12, 17).Value = 48.15661950747901
$ws.Cells.Item(12, 18).Value = 433.4095755673111
$ws.Cells.Item(12, 19).Value = 0.0920368135304723
$ws.Cells.Item(12, 20).Value = 0.0920368135304723

$ws.Cells.Item(13, 7).Value = 2.651481
$ws.Cells.Item(13, 8).Value = 7.954443000000001
$ws.Cells.Item(13, 9).Value = 0.3159313972167653
$ws.Cells.Item(13, 10).Value = 0.3159313972167653
$ws.Cells.Item(13, 13).Value = 0.3271573333333334
$ws.Cells.Item(13, 14).Value = 0.981472
$ws.Cells.Item(13, 15).Value = 0.005247566758452071
$ws.Cells.Item(13, 16).Value = 0.00524756675845207
$ws.Cells.Item(13, 17).Value = 0.8674514533440002
$ws.Cells.Item(13, 18).Value = 7.807063080096001
$ws.Cells.Item(13, 19).Value = 0.001657871097986014
$ws.Cells.Item(13, 20).Value = 0.001657871097986014

$ws.Cells.Item(14, 7).Value = 3.633574000000001
$ws.Cells.Item(14, 8).Value = 10.900722
$ws.Cells.Item(14, 9).Value = 0.4329505324422505
$ws.Cells.Item(14, 10).Value = 0.4329505324422504
$ws.Cells.Item(14, 13).Value = 10.43365533333333
$ws.Cells.Item(14, 14).Value = 31.300966
$ws.Cells.Item(14, 15).Value = 0.1673546557507891
$ws.Cells.Item(14, 16).Value = 0.1673546557507891
$ws.Cells.Item(14, 17).Value = 37.91145874416134
$ws.Cells.Item(14, 18).Value = 341.2031286974521
$ws.Cells.Item(14, 19).Value = 0.07245628731399367
$ws.Cells.Item(14, 20).Value = 0.07245628731399366

$ws.Cells.Item(15, 7).Value = 3.633574000000001
$ws.Cells.Item(15, 8).Value = 10.900722
$ws.Cells.Item(15, 9).Value = 0.4329505324422505
$ws.Cells.Item(15, 10).Value = 0.4329505324422504
$ws.Cells.Item(15, 15).Value = 0.5360787938719054
$ws.Cells.Item(15, 16).Value = 0.5360787938719054
$ws.Cells.Item(15, 17).Value = 121.439878599844
$ws.Cells.Item(15, 18).Value = 1092.958907398596
$ws.Cells.Item(15, 19).Value = 0.2320955992378409
$ws.Cells.Item(15, 20).Value = 0.2320955992378408

$ws.Cells.Item(16, 7).Value = 3.633574000000001
$ws.Cells.Item(16, 8).Value = 10.900722
$ws.Cells.Item(16, 9).Value = 0.4329505324422505
$ws.Cells.Item(16, 10).Value = 0.4329505324422504
$ws.Cells.Item(16, 13).Value = 18.162159
$ws.Cells.Item(16, 14).Value = 54.486477
$ws.Cells.Item(16, 15).Value = 0.2913189836188534
$ws.Cells.Item(16, 16).Value = 0.2913189836188534
$ws.Cells.Item(16, 17).Value = 65.99354872626601
$ws.Cells.Item(16, 18).Value = 593.9419385363941
$ws.Cells.Item(16, 19).Value = 0.1261267090683178
$ws.Cells.Item(16, 20).Value = 0.1261267090683178

$ws.Cells.Item(17, 7).Value = 3.633574000000001
$ws.Cells.Item(17, 8).Value = 10.900722
$ws.Cells.Item(17, 9).Value = 0.4329505324422505
$ws.Cells.Item(17, 10).Value = 0.4329505324422504
$ws.Cells.Item(17, 13).Value = 0.3271573333333334
$ws.Cells.Item(17, 14).Value = 0.981472
$ws.Cells.Item(17, 15).Value = 0.005247566758452071
$ws.Cells.Item(17, 16).Value = 0.00524756675845207
$ws.Cells.Item(17, 17).Value = 1.188750380309334
$ws.Cells.Item(17, 18).Value = 10.698753422784
$ws.Cells.Item(17, 19).Value = 0.002271936822098078
$ws.Cells.Item(17, 20).Value = 0.002271936822098078
